$wb = $excel.ActiveWorkbook

# Rename the existing "strategy_id-5008" sheet to "strategy_id-5007"
$ws3 = $wb.Worksheets.Item("strategy_id-5008")
$ws3.Name = "strategy_id-5007"

# Add a new sheet "strategy_id-5009" right after it, as a copy of ws3's content
$ws3.Copy([System.Reflection.Missing]::Value, $ws3)
$ws4 = $wb.Worksheets.Item($ws3.Index + 1)
$ws4.Name = "strategy_id-5009"

# Restore original active sheet/tab selection so the copy/rename operations
# don't leave the new sheet as the active one.
$wb.Worksheets.Item(1).Activate()
